$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text/string updates (safe from numeric auto-coercion) ---
$ws.Range("D2").Value = "64.910.23"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "3.174.14"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.174.58"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "3.694.68"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "64.904.65"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "3.171.63"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -3.48%  "
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").Value = "  -6.02%  "
$ws.Range("E30").Value = "  -5.69%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("D36").Value = "0.0₃0795"
$ws.Range("E36").Value = "  +5.91%  "
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("E40").Value = "  +3.80%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").Value = "2.863.74"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("E47").Value = "  +6.35%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  -1.01%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("E49").Value = "  +10.71%  "

# --- Cells whose new text looks like a plain decimal number; Excel would
#     auto-convert these to numeric via COM Value assignment. Force text by
#     temporarily setting NumberFormat to Text, then restore the default
#     "Normal" style afterwards so no visible formatting change remains. ---
$textForcedCells = @("D5", "D6", "D11", "D14", "D20", "D21", "D24", "D25", "D27", "D28", "D30", "D31", "D34", "D37", "D38", "D40", "D41", "D47", "D48", "D49")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "617.36"
$ws.Range("D6").Value = "147.16"
$ws.Range("D11").Value = "5.51"
$ws.Range("D14").Value = "35.96"
$ws.Range("D20").Value = "481.55"
$ws.Range("D21").Value = "14.77"
$ws.Range("D24").Value = "13.87"
$ws.Range("D25").Value = "84.74"
$ws.Range("D27").Value = "2.84"
$ws.Range("D28").Value = "8.70"
$ws.Range("D30").Value = "2.12"
$ws.Range("D31").Value = "6.96"
$ws.Range("D34").Value = "26.77"
$ws.Range("D37").Value = "6.06"
$ws.Range("D38").Value = "3.22"
$ws.Range("D40").Value = "467.71"
$ws.Range("D41").Value = "0.0403"
$ws.Range("D47").Value = "2.46"
$ws.Range("D48").Value = "26.86"
$ws.Range("D49").Value = "36.90"

foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
